$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted before the current row 19, shifting
# all subsequent rows (old 19..42) down by one (new 20..43).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record's data.
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44792
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112001
$ws.Cells.Item(19, 7).Value = "Berenjena"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 60
$ws.Cells.Item(19, 11).Value = 12000
$ws.Cells.Item(19, 12).Value = 13000
$ws.Cells.Item(19, 13).Value = 12500
$ws.Cells.Item(19, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(19, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(19, 16).Value = 208
$ws.Cells.Item(19, 17).Value = 60
$ws.Cells.Item(19, 18).Value = "Hortaliza"
